# Text updates as supplied by PM&C.
#
# The "Description" sheet's source note (old B9: "Sourced from ABS Australian
# Health Survey") is replaced by a "Source" label in A9/B9 plus three
# individual citation lines appended below it (B9:B11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")
$ws.Activate() | Out-Null

# B9 currently holds the only use of the old "Sourced from ABS Australian
# Health Survey" shared string - overwrite it with "Source" first so that
# string slot gets reused/renamed rather than orphaned, matching A9 which
# becomes the visible "Source" label cell.
$ws.Range("B9").Value = "Source"
$ws.Range("A9").Value = "Source"

# Give B9 the same (unwrapped, borderless) formatting as A9 before putting
# the real citation text into it.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = "ABS (unpublished) Australian Health Survey 2014" + [char]0x2013 + "15"

# New citation rows.
$ws.Range("B10").Value = "ABS (unpublished) Australian Health Survey, 2011-13 (2011-12 Core component)"
$ws.Range("B11").Value = "ABS (unpublished), National Health Survey 2007-08."

$citations = $ws.Range("B10:B11")
$citations.RowHeight = 15
$citations.Font.Color = 0
$citations.Font.Name = "Arial"
$citations.Font.Size = 12

$ws.Range("B12").Select() | Out-Null
